$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "formulae"
$ws.Range("G2").Formula = "=SUM(D2,F2)"
$ws.Range("G3").Formula = "=SUM(D3,F3)"
